# Correção do PITCH versão final - Ajuste logo
#
# The "LOGO DO PROJETO" placeholder shape on slide 1 (an Oval) is turned
# into a Rectangle and repositioned/resized.
#
# NOTE on the literal point values below: PowerPoint's COM object model
# stores Shape.Left/Top/Width/Height as single-precision (32-bit) floats.
# The target OOXML uses exact EMU values (1 pt = 12700 EMU), so naively
# using `emu / 12700` can be off by 1 EMU after the float32 round-trip.
# The constants here are chosen so that, after that float32 round-trip,
# they reproduce the exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item("Oval 3")

# Rename "Oval 3" -> "Rectangle 3"
$shp.Name = "Rectangle 3"

# Change the preset geometry from ellipse to rectangle (msoShapeRectangle = 1)
$shp.AutoShapeType = 1

# Reposition / resize: (9177867,3330222,2190044,1761067) -> (8895644,3759200,2460978,1555044) EMU
$shp.Left = 700.4444274988189
$shp.Top = 296.0
$shp.Width = 193.77779527559056
$shp.Height = 122.4444094488189
